$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10Folds")

# --- Correct the raw per-fold "time" measurements (A21:D30) ---
# The original values were recorded in the wrong (inflated) unit; replace
# with the corrected measurements.
$timeData = @(
    @(6.1446211338043204, 2000.1803841590799, 228.58463430404601, 10.0875227451324),
    @(5.94313168525695,   2532.5105450153301, 198.90165853500301, 10.455570697784401),
    @(5.9491238594055096, 2707.3796503543799, 194.60154581069901, 10.4448153972625),
    @(5.97568488121032,   1974.6941828727699, 187.938019275665,   10.471176624298),
    @(6.0945928096771196, 2567.4352707862799, 190.35401844978301, 10.433064699172901),
    @(6.0776596069335902, 2518.2201952934201, 196.90826129913299, 10.4409625530242),
    @(6.4992554187774596, 2282.8988213539101, 182.76189756393401, 10.433826208114599),
    @(7.9892451763152996, 2302.8733179569199, 171.421284675598,   10.448426485061599),
    @(6.45373106002807,   2343.7110028266902, 182.22777366638101, 10.4903218746185),
    @(7.6102759838104204, 2790.04174423217,   182.77237868309001, 10.420721292495699)
)

$row = 21
foreach ($vals in $timeData) {
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $row++
}

# --- Corresponding ANOVA SUMMARY block (I21:K24), recomputed from the new data ---
$ws.Range("I21").Value = 64.737321615219059
$ws.Range("J21").Value = 6.4737321615219061
$ws.Range("K21").Value = 0.53369989482009339

$ws.Range("I22").Value = 24019.945114850951
$ws.Range("J22").Value = 2401.9945114850952
$ws.Range("K22").Value = 74800.016983884911

$ws.Range("I23").Value = 1916.4714722633319
$ws.Range("J23").Value = 191.64714722633317
$ws.Range("K23").Value = 236.04030920619621

$ws.Range("I24").Value = 104.12640857696481
$ws.Range("J24").Value = 10.412640857696481
$ws.Range("K24").Value = 0.013450523444093844

# --- ANOVA table (H29:L32), recomputed from the new data ---
$ws.Range("H29").Value = 41027424.434293725
$ws.Range("J29").Value = 13675808.144764574
$ws.Range("K29").Value = 729.02062912827751
$ws.Range("L29").Value = 0.000000000000000000000000000000028440341794522369

$ws.Range("H30").Value = 675329.43999160163
$ws.Range("J30").Value = 18759.151110877821

$ws.Range("H32").Value = 41702753.874285325

# --- View state: scroll position and selection ---
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B32:B34").Select() | Out-Null

# --- Column K (11) picks up a best-fit width once its numbers are shorter ---
$ws.Columns.Item(11).EntireColumn.AutoFit()

# --- Conditional-formatting: bump the accuracy block's priority out of the
#     way, and extend the time block up to include the new header row 20 ---
$ws.Range("H2:M16").FormatConditions.Item(1).Priority = 3
$ws.Range("H18:M32").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H20:M32"))
